$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 6374.75
$ws.Range("I47").Value = 6374.75
$ws.Range("K47").Value = 6374.75
$ws.Range("M47").Value = -5402.75
$ws.Range("H54").Value = 9980
$ws.Range("I54").Value = 9980
$ws.Range("K54").Value = 9980
$ws.Range("M54").Value = -9494
$ws.Range("H80").Value = 609.25
$ws.Range("I80").Value = 512.3333
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 1536.9999
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -538.9999
$ws.Range("N80").Value = -4696
$ws.Range("H83").Value = 609.25
$ws.Range("I83").Value = 512.3333
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 4610.9997
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = 381.0002999999997
$ws.Range("N83").Value = -18084
$ws.Range("H100").Value = 2437
$ws.Range("I100").Value = 2452.5
$ws.Range("J100").Value = 2406
$ws.Range("K100").Value = 2452.5
$ws.Range("L100").Value = 2406
$ws.Range("M100").Value = -1911.5
$ws.Range("N100").Value = -3488
$ws.Range("H113").Value = 3998
$ws.Range("J113").Value = 3999.5
$ws.Range("L113").Value = 3999.5
$ws.Range("N113").Value = -10507.5
$ws.Range("H131").Value = 3556.5217
$ws.Range("I131").Value = 1108.8
$ws.Range("K131").Value = 3326.4
$ws.Range("M131").Value = 1713.6
$ws.Range("H132").Value = 1687
$ws.Range("I132").Value = 1594.48
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4783.440000000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2253.440000000001
$ws.Range("N132").Value = -17060
$ws.Range("H137").Value = 1837.6522
$ws.Range("I137").Value = 928.61536
$ws.Range("J137").Value = 3019.4
$ws.Range("K137").Value = 2785.84608
$ws.Range("L137").Value = 9058.200000000001
$ws.Range("M137").Value = -235.8460800000003
$ws.Range("N137").Value = -14158.2

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1756
$ws.Range("I45").Value = 1756
$ws.Range("K45").Value = 1756
$ws.Range("M45").Value = -1379
$ws.Range("H61").Value = 2212.2144
$ws.Range("I61").Value = 2212.2144
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2212.2144
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2000.2144
$ws.Range("H63").Value = 5790.3125
$ws.Range("I63").Value = 4678.2856
$ws.Range("K63").Value = 4678.2856
$ws.Range("M63").Value = -3992.2856
$ws.Range("H66").Value = 5790.3125
$ws.Range("I66").Value = 4678.2856
$ws.Range("K66").Value = 23391.428
$ws.Range("M66").Value = -19959.428
$ws.Range("H74").Value = 2041.2
$ws.Range("I74").Value = 987.8461
$ws.Range("K74").Value = 987.8461
$ws.Range("M74").Value = -113.8461
$ws.Range("H77").Value = 2041.2
$ws.Range("I77").Value = 987.8461
$ws.Range("K77").Value = 4939.2305
$ws.Range("M77").Value = -571.2304999999997
$ws.Range("H110").Value = 5959.3335
$ws.Range("I110").Value = 6921.2
$ws.Range("J110").Value = 1150
$ws.Range("K110").Value = 6921.2
$ws.Range("L110").Value = 1150
$ws.Range("M110").Value = -4876.2
$ws.Range("N110").Value = -5240
$ws.Range("H122").Value = 386930.78
$ws.Range("I122").Value = 527644.25
$ws.Range("J122").Value = 4994.143
$ws.Range("K122").Value = 1582932.75
$ws.Range("L122").Value = 14982.429
$ws.Range("M122").Value = -1580482.75
$ws.Range("N122").Value = -19882.429
$ws.Range("H136").Value = 2212.2144
$ws.Range("I136").Value = 2212.2144
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6636.6432
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -4086.6432

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2565
$ws.Range("I94").Value = 2385.625
$ws.Range("K94").Value = 2385.625
$ws.Range("M94").Value = -1934.625
$ws.Range("H99").Value = 4027.0908
$ws.Range("I99").Value = 3829.8
$ws.Range("K99").Value = 3829.8
$ws.Range("M99").Value = -2331.8
$ws.Range("H105").Value = 4128.8335
$ws.Range("I105").Value = 3254.8948
$ws.Range("K105").Value = 3254.8948
$ws.Range("M105").Value = -1507.8948
$ws.Range("H107").Value = 631.375
$ws.Range("I107").Value = 623.9545000000001
$ws.Range("K107").Value = 623.9545000000001
$ws.Range("M107").Value = 1296.0455

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7108.875
$ws.Range("I31").Value = 3497.5
$ws.Range("J31").Value = 8312.666999999999
$ws.Range("K31").Value = 3497.5
$ws.Range("L31").Value = 8312.666999999999
$ws.Range("M31").Value = -3202.5
$ws.Range("N31").Value = -8902.666999999999
$ws.Range("H34").Value = 7108.875
$ws.Range("I34").Value = 3497.5
$ws.Range("J34").Value = 8312.666999999999
$ws.Range("K34").Value = 3497.5
$ws.Range("L34").Value = 8312.666999999999
$ws.Range("M34").Value = -3295.5
$ws.Range("N34").Value = -8716.666999999999
$ws.Range("H58").Value = 2980.2083
$ws.Range("J58").Value = 3537.6667
$ws.Range("L58").Value = 3537.6667
$ws.Range("N58").Value = -3943.6667
$ws.Range("H99").Value = 9962.457
$ws.Range("I99").Value = 7074.7144
$ws.Range("J99").Value = 11887.619
$ws.Range("K99").Value = 7074.7144
$ws.Range("L99").Value = 11887.619
$ws.Range("M99").Value = -5576.7144
$ws.Range("N99").Value = -14883.619
$ws.Range("H126").Value = 9962.457
$ws.Range("I126").Value = 7074.7144
$ws.Range("J126").Value = 11887.619
$ws.Range("K126").Value = 21224.1432
$ws.Range("L126").Value = 35662.857
$ws.Range("M126").Value = -18754.1432
$ws.Range("N126").Value = -40602.857
$ws.Range("H132").Value = 3218.6365
$ws.Range("I132").Value = 1929.3572
$ws.Range("K132").Value = 5788.071599999999
$ws.Range("M132").Value = -3258.071599999999
$ws.Range("H136").Value = 2980.2083
$ws.Range("J136").Value = 3537.6667
$ws.Range("L136").Value = 10613.0001
$ws.Range("N136").Value = -15713.0001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 46432548
$ws.Range("I4").Value = 51075228
$ws.Range("J4").Value = 5749.5
$ws.Range("K4").Value = 153225684
$ws.Range("L4").Value = 17248.5
$ws.Range("M4").Value = -153225572
$ws.Range("N4").Value = -17472.5
$ws.Range("H50").Value = 187.25
$ws.Range("I50").Value = 118.8
$ws.Range("J50").Value = 301.33334
$ws.Range("K50").Value = 356.4
$ws.Range("L50").Value = 904.0000200000001
$ws.Range("M50").Value = 124.6
$ws.Range("N50").Value = -1866.00002
$ws.Range("H53").Value = 187.25
$ws.Range("I53").Value = 118.8
$ws.Range("J53").Value = 301.33334
$ws.Range("K53").Value = 356.4
$ws.Range("L53").Value = 904.0000200000001
$ws.Range("M53").Value = 124.6
$ws.Range("N53").Value = -1866.00002
$ws.Range("H137").Value = 10122
$ws.Range("J137").Value = 10433
$ws.Range("L137").Value = 31299
$ws.Range("N137").Value = -41499

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4345735.5
$ws.Range("I29").Value = 12286167
$ws.Range("K29").Value = 12286167
$ws.Range("M29").Value = -12285877
$ws.Range("H122").Value = 38671.215
$ws.Range("I122").Value = 2391.6
$ws.Range("J122").Value = 129370.25
$ws.Range("K122").Value = 7174.799999999999
$ws.Range("L122").Value = 388110.75
$ws.Range("M122").Value = -4724.799999999999
$ws.Range("N122").Value = -393010.75
$ws.Range("H132").Value = 1922.56
$ws.Range("I132").Value = 1168.85
$ws.Range("K132").Value = 3506.55
$ws.Range("M132").Value = -976.5499999999997

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1373.4
$ws.Range("I16").Value = 1373.4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1373.4
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1203.4
$ws.Range("H40").Value = 1189
$ws.Range("I40").Value = 1053
$ws.Range("K40").Value = 1053
$ws.Range("M40").Value = -917
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H68").Value = 1994.4
$ws.Range("I68").Value = 1968.625
$ws.Range("K68").Value = 1968.625
$ws.Range("M68").Value = -1219.625
$ws.Range("H71").Value = 1994.4
$ws.Range("I71").Value = 1968.625
$ws.Range("K71").Value = 9843.125
$ws.Range("M71").Value = -6099.125
$ws.Range("H93").Value = 200
$ws.Range("J93").Value = 298
$ws.Range("L93").Value = 298
$ws.Range("N93").Value = -2794
$ws.Range("H122").Value = 6153.385
$ws.Range("I122").Value = 3997.25
$ws.Range("K122").Value = 11991.75
$ws.Range("M122").Value = -9541.75
$ws.Range("H132").Value = 4212.676
$ws.Range("I132").Value = 3629.3914
$ws.Range("K132").Value = 10888.1742
$ws.Range("M132").Value = -8358.174199999999
$ws.Range("H136").Value = 3029.7273
$ws.Range("I136").Value = 3158.8
$ws.Range("J136").Value = 1739
$ws.Range("K136").Value = 9476.400000000001
$ws.Range("L136").Value = 5217
$ws.Range("M136").Value = -6926.400000000001
$ws.Range("N136").Value = -10317

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8352.117
$ws.Range("J62").Value = 8352.117
$ws.Range("L62").Value = 8352.117
$ws.Range("N62").Value = -9600.117
$ws.Range("H65").Value = 8352.117
$ws.Range("J65").Value = 8352.117
$ws.Range("L65").Value = 41760.585
$ws.Range("N65").Value = -48000.585
$ws.Range("H122").Value = 3088.75
$ws.Range("I122").Value = 1102.8572
$ws.Range("K122").Value = 3308.5716
$ws.Range("M122").Value = -858.5715999999998
$ws.Range("H132").Value = 1018.6667
$ws.Range("I132").Value = 1022.29034
$ws.Range("K132").Value = 3066.87102
$ws.Range("M132").Value = -536.87102
